$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.531.09'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').Value = '2.432.60'
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('D5').Value = "'513.94"
$ws.Range('E5').Value = '  -2.57%  '
$ws.Range('D6').Value = "'129.58"
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').Value = '2.446.28'
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').Value = "'0.0951"
$ws.Range('E11').Value = '  -5.43%  '
$ws.Range('D12').Value = "'5.18"
$ws.Range('E12').Value = '  -4.09%  '
$ws.Range('D13').Value = "'0.332"
$ws.Range('E13').Value = '  -3.36%  '
$ws.Range('D14').Value = '2.867.39'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '57.462.77'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').Value = "'21.81"
$ws.Range('E16').Value = '  -3.25%  '
$ws.Range('E17').Value = '  -3.69%  '
$ws.Range('D18').Value = '2.441.07'
$ws.Range('E18').Value = '  -1.89%  '
$ws.Range('D19').Value = "'10.44"
$ws.Range('E19').Value = '  -4.88%  '
$ws.Range('D20').Value = "'315.42"
$ws.Range('E21').Value = '  -2.75%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = "'5.67"
$ws.Range('E23').Value = '  -3.02%  '
$ws.Range('D24').Value = "'63.42"
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('D25').Value = "'0.407"
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').Value = "'170.18"
$ws.Range('E29').Value = '  +2.54%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = "'6.25"
$ws.Range('E30').Value = '  -2.98%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0720'
$ws.Range('E31').Value = '  -4.71%  '
$ws.Range('E32').Value = '  -3.36%  '
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = "'17.66"
$ws.Range('E36').Value = '  -3.50%  '
$ws.Range('D37').Value = "'1.28"
$ws.Range('E37').Value = '  -4.79%  '
$ws.Range('D38').Value = "'3.90"
$ws.Range('E38').Value = '  -2.53%  '
$ws.Range('D39').Value = "'36.19"
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').Value = "'1.45"
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('D41').Value = "'0.777"
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'3.37"
$ws.Range('E42').Value = '  -4.75%  '
$ws.Range('D43').Value = "'268.01"
$ws.Range('E43').Value = '  -3.17%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = "'4.90"
$ws.Range('E44').Value = '  -1.37%  '
$ws.Range('D45').Value = "'0.586"
$ws.Range('E45').Value = '  -1.92%  '
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').Value = "'120.44"
$ws.Range('E47').Value = '  -5.65%  '
$ws.Range('D48').Value = "'0.0484"
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').Value = "'0.0210"
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('D50').Value = "'16.53"
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('D51').Value = '1.707.61'
$ws.Range('E51').Value = '  -2.09%  '

# Reset style on forced-text numeric-looking cells to avoid leftover quotePrefix style
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
